$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44160
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 175
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 18743
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1250
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44179
$ws.Range("K3").Value = 'Modesto'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 22000
$ws.Range("Q3").Value = '$/caja 15 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1467
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44179
$ws.Range("K4").Value = 'Modesto'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1111
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44188
$ws.Range("K5").Value = 'Modesto'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 35
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = '$/bandeja 18 kilos'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 1111
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 44159
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 85
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = '$/caja 15 kilos'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 1333
$ws.Range("T6").Value = 15

$ws.Range("D7").Value = 44186
$ws.Range("K7").Value = 'Modesto'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = '$/bandeja 18 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44175
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = '$/bandeja 18 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1111
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44175
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 55
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("Q9").Value = '$/bandeja 18 kilos'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 44175
$ws.Range("K10").Value = 'Castle Brite'
$ws.Range("L10").Value = 'Tercera'
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range("Q10").Value = '$/bandeja 18 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 778
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44168
$ws.Range("K11").Value = 'Castle Brite'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 450
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 22444
$ws.Range("Q11").Value = '$/bandeja 18 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1247
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44217
$ws.Range("K12").Value = 'Modesto'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/bandeja 18 kilos'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 44162
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 75
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 18933
$ws.Range("Q13").Value = '$/caja 15 kilos'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 1262
$ws.Range("T13").Value = 15

$ws.Range("D14").Value = 44162
$ws.Range("K14").Value = 'Castle Brite'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 55
$ws.Range("N14").Value = 23000
$ws.Range("O14").Value = 23000
$ws.Range("P14").Value = 23000
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1278
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 44202
$ws.Range("K15").Value = 'Modesto'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 18000
$ws.Range("Q15").Value = '$/caja 15 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1200
$ws.Range("T15").Value = 15

$ws.Range("D16").Value = 44167
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("Q16").Value = '$/caja 15 kilos'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1333
$ws.Range("T16").Value = 15

$ws.Range("D17").Value = 44174
$ws.Range("K17").Value = 'Castle Brite'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21083
$ws.Range("Q17").Value = '$/bandeja 18 kilos'
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1171
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 44172
$ws.Range("K18").Value = 'Castle Brite'
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 22000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 22000
$ws.Range("Q18").Value = '$/bandeja 18 kilos'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1222
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 44172
$ws.Range("K19").Value = 'Castle Brite'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 65
$ws.Range("N19").Value = 23000
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 23000
$ws.Range("Q19").Value = '$/bandeja 18 kilos'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1278
$ws.Range("T19").Value = 18

$ws.Range("D20").Value = 44176
$ws.Range("K20").Value = 'Modesto'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = '$/caja 18 kilos'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 1111
$ws.Range("T20").Value = 18

$ws.Range("D21").Value = 44166
$ws.Range("K21").Value = 'Castle Brite'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 400
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 17750
$ws.Range("Q21").Value = '$/caja 15 kilos'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 1183
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 44166
$ws.Range("K22").Value = 'Castle Brite'
$ws.Range("L22").Value = 'Segunda'
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("Q22").Value = '$/caja 15 kilos'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 800
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 44181
$ws.Range("K23").Value = 'Modesto'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 140
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 17000
$ws.Range("P23").Value = 17000
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 944
$ws.Range("T23").Value = 18

$ws.Range("D24").Value = 44201
$ws.Range("K24").Value = 'Modesto'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 45
$ws.Range("N24").Value = 18000
$ws.Range("O24").Value = 18000
$ws.Range("P24").Value = 18000
$ws.Range("Q24").Value = '$/caja 15 kilos'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 1200
$ws.Range("T24").Value = 15
